# Insert a new data row before the current row 295, shifting the existing
# rows 295-356 down to 296-357, then populate the new row with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(295).Insert()

$ws.Cells.Item(295, 1).Value2  = 5
$ws.Cells.Item(295, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(295, 3).Value2  = "Maule"
$ws.Cells.Item(295, 4).Value2  = 44711
$ws.Cells.Item(295, 5).Value2  = 7
$ws.Cells.Item(295, 6).Value2  = 100112032
$ws.Cells.Item(295, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(295, 8).Value2  = "Sin especificar"
$ws.Cells.Item(295, 9).Value2  = "Primera"
$ws.Cells.Item(295, 10).Value2 = 300
$ws.Cells.Item(295, 11).Value2 = 16000
$ws.Cells.Item(295, 12).Value2 = 16000
$ws.Cells.Item(295, 13).Value2 = 16000
$ws.Cells.Item(295, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(295, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(295, 16).Value2 = 320
$ws.Cells.Item(295, 17).Value2 = 50
$ws.Cells.Item(295, 18).Value2 = "Hortaliza"
